# Reformat the "login" sheet: repeat the username/password sample data
# across a 4-column x 7-row grid (header row + 6 data rows) and include
# a proper header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "password"

# Data rows 2-7: same username/password repeated across all four columns
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "ideam363@gmail.com"
    $ws.Cells.Item($r, 2).Value = "Mandar@9766"
    $ws.Cells.Item($r, 3).Value = "Mandar@9766"
    $ws.Cells.Item($r, 4).Value = "Mandar@9766"
}

# Match the selection recorded in the edited workbook
$ws.Range("A2:D7").Select() | Out-Null
